# Append the 2025-10-12 allocation row (Date, BTC, KAS) to Sheet1, as
# produced by the daily "profit files" run on 2025-10-12.
#
# New row 41: A41="10/12/2025", B41=0.1777389139054508, C41=0.8222610860945492

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Assigning the literal string "10/12/2025" straight to a cell's
# .Value makes Excel's normal type-inference kick in and silently store it
# as a date serial number (45942) instead of text - which also forces a new
# number-format style onto the cell. The existing Date column (A2:A40) is
# plain text, so we need to land "10/12/2025" as literal text too.
#
# Route it through a scratch cell using =TEXT(...,"@") (which evaluates to
# a genuine text value, not a date) and then Copy/paste it as a value into
# A41 - this stores a plain text value with no extra formatting, matching
# the rest of the column. The scratch cell is cleared afterwards so it
# leaves no trace in the saved workbook.
$scratch = $ws.Cells.Item(100, 26)
$scratch.Formula = "=TEXT(""10/12/2025"",""@"")"
$scratch.Copy($ws.Cells.Item(41, 1))
$scratch.Clear()

$ws.Cells.Item(41, 2).Value = 0.1777389139054508
$ws.Cells.Item(41, 3).Value = 0.8222610860945492
